$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo in header E1: "chi2p-value" -> "chi2 p-value"
$ws.Range("E1").Value = "chi2 p-value"

# Update row 2 values (new model fit indices)
$ws.Range("D2").Value = 48.655132999999999
$ws.Range("E2").Value = [double]"2.6102179999999999E-9"
$ws.Range("F2").Value = 98.099002999999996
$ws.Range("G2").Value = 0.51547600000000005
$ws.Range("H2").Value = 0.50402000000000002
$ws.Range("I2").Value = 0.206432
$ws.Range("J2").Value = 0.50402000000000002
$ws.Range("K2").Value = 0.22476199999999999
$ws.Range("L2").Value = 0.32831500000000002
$ws.Range("M2").Value = 8.8132889999999993
$ws.Range("N2").Value = 20.846886000000001
$ws.Range("O2").Value = 0.59335499999999997

# Apply scientific-notation number format to E2 (numFmtId 11 -> "0.00E+00")
$ws.Range("E2").NumberFormat = "0.00E+00"

# Update row 10 values (new model fit indices)
$ws.Range("D10").Value = 7.262365
$ws.Range("E10").Value = [double]"6.3989000000000004E-2"
$ws.Range("F10").Value = 100.60639
$ws.Range("G10").Value = 0.95446500000000001
$ws.Range("H10").Value = 0.92781400000000003
$ws.Range("I10").Value = 0.83156600000000003
$ws.Range("J10").Value = 0.92781400000000003
$ws.Range("K10").Value = 0.89375199999999999
$ws.Range("L10").Value = 0.132441
$ws.Range("M10").Value = 13.822869000000001
$ws.Range("N10").Value = 30.669903999999999
$ws.Range("O10").Value = [double]"8.8565000000000005E-2"

# Update the selection to match the new active cell/range
$ws.Range("I2:O2").Select()
